{"js": "const body = context.document.body;\nconst results1 = body.search(\"sm\u011bnka je splatna p\u0159i p\u0159edlo\u017een\u00ed,\", { matchCase: true });\nresults1.load(\"items\");\nawait context.sync();\nfor (const r of results1.items) {\n  r.insertText(\"sm\u011bnka je splatn\u00e1 p\u0159i p\u0159edlo\u017een\u00ed,\", Word.InsertLocation.replace);\n}\n\nconst results2 = body.search(\"bod sedn\u00fd je rovn\u011b\u017e spln\u011bn (podepsal).\", { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\nfor (const r of results2.items) {\n  r.insertText(\"bod sedm\u00fd je rovn\u011b\u017e spln\u011bn (podepsal).\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"sm\u011bnka je splatna p\u0159i p\u0159edlo\u017een\u00ed,\"\n$find.Replacement.Text = \"sm\u011bnka je splatn\u00e1 p\u0159i p\u0159edlo\u017een\u00ed,\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n$find2 = $d.Content.Find\n$find2.Text = \"bod sedn\u00fd je rovn\u011b\u017e spln\u011bn (podepsal).\"\n$find2.Replacement.Text = \"bod sedm\u00fd je rovn\u011b\u017e spln\u011bn (podepsal).\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
